$wb = $excel.ActiveWorkbook

$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb102fc023d36a547c393785c62d8f12d0bcdd30/e2e/7970f794-df68-4d2b-a9f2-4fdecfd3e2e9.md"
$currentUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/93e44a2e9535c16c64ae5ba9505a23d6a0c88470/e2e/7970f794-df68-4d2b-a9f2-4fdecfd3e2e9.md"
$errorDetail = "The version of handback file is not the latest, current: $currentUrl, latest: $latestUrl."
$targetFileDisplay = "7970f794-df68-4d2b-a9f2-4fdecfd3e2e9.md"

function Update-LangSheet {
    param(
        [string]$sheetName,
        [string]$handbackDateTime
    )

    $ws = $wb.Worksheets.Item($sheetName)

    # Error Detail column (P) needs to be wide enough to show the message.
    $ws.Columns.Item(16).ColumnWidth = 40

    # Row 7 ("7970f794-...") gets a handback record: target file (I7), handback
    # file (J7, mirrors the latest handoff xliff in G7), handback datetime (K7)
    # and an error detail (P7) noting the handback file is stale.
    $g7 = $ws.Range("G7").Value2

    $ws.Range("J7").Value = $g7
    $ws.Range("K7").Value = $handbackDateTime
    $ws.Range("P7").Value = $errorDetail

    $ws.Range("I7").Value = $targetFileDisplay
    $ws.Range("I7").Style = "Hyperlink"
    $ws.Hyperlinks.Add($ws.Range("I7"), $latestUrl, "", "", $targetFileDisplay) | Out-Null
}

Update-LangSheet "zh-cn" "2016-08-20 18:53:13"
Update-LangSheet "de-de" "2016-08-20 18:53:19"
